# Vent process added before calibration process start - per commit message.
# Update the Sheet1 configuration values and the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# G7: 0.01 -> 0
$ws.Range("G7").Value = 0

# H8: 25.6 -> 28 (Temp[c] updated)
$ws.Range("H8").Value = 28

# K8: 1 -> 5 (Samples amount updated)
$ws.Range("K8").Value = 5

# H9: clear the Temp[c] value (cell becomes empty)
$ws.Range("H9").ClearContents()

# K11: 0.2 -> 0.1 (Max time wait to temp stable updated)
$ws.Range("K11").Value = 0.1

# Move / leave the active selection at K12, matching the final saved view state.
$ws.Range("K12").Select()
